$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21, column Y (25): change date number format to include time (s=3 -> s=2)
$ws.Range("Y21").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append new row 22, a copy of the former row 21 data, with Y22 keeping the date-only format
$ws.Cells.Item(22,1).Value = 565
$ws.Cells.Item(22,2).Value = 474
$ws.Cells.Item(22,3).Value = 450
$ws.Cells.Item(22,4).Value = 522
$ws.Cells.Item(22,5).Value = 515
$ws.Cells.Item(22,6).Value = 520
$ws.Cells.Item(22,7).Value = 474
$ws.Cells.Item(22,8).Value = 570
$ws.Cells.Item(22,9).Value = 490
$ws.Cells.Item(22,10).Value = 450
$ws.Cells.Item(22,11).Value = 571
$ws.Cells.Item(22,12).Value = 480
$ws.Cells.Item(22,13).Value = 485
$ws.Cells.Item(22,14).Value = 505
$ws.Cells.Item(22,15).Value = 545
$ws.Cells.Item(22,16).Value = 480
$ws.Cells.Item(22,17).Value = 618
$ws.Cells.Item(22,18).Value = 490
$ws.Cells.Item(22,19).Value = 474
$ws.Cells.Item(22,20).Value = 480
$ws.Cells.Item(22,21).Value = 619
$ws.Cells.Item(22,22).Value = 550
$ws.Cells.Item(22,23).Value = 599
$ws.Cells.Item(22,24).Value = 495
$ws.Cells.Item(22,25).Value = 45754
$ws.Cells.Item(22,26).Value = 850
$ws.Cells.Item(22,27).Value = 555
$ws.Cells.Item(22,28).Value = 543.5
$ws.Cells.Item(22,29).Value = 500
$ws.Cells.Item(22,30).Value = 545
$ws.Cells.Item(22,31).Value = 507
$ws.Cells.Item(22,32).Value = 509
$ws.Cells.Item(22,33).Value = 745
$ws.Cells.Item(22,34).Value = 473
$ws.Cells.Item(22,35).Value = 735
$ws.Cells.Item(22,36).Value = 474
$ws.Cells.Item(22,37).Value = 488
$ws.Cells.Item(22,38).Value = 570
$ws.Cells.Item(22,39).Value = 555
$ws.Cells.Item(22,40).Value = 488
$ws.Cells.Item(22,41).Value = 535
$ws.Cells.Item(22,42).Value = 547
$ws.Cells.Item(22,43).Value = 568
$ws.Cells.Item(22,44).Value = 547
$ws.Cells.Item(22,45).Value = 645
$ws.Cells.Item(22,46).Value = 634
$ws.Cells.Item(22,47).Value = 496
$ws.Cells.Item(22,48).Value = 485

$ws.Range("Y22").NumberFormat = "YYYY-MM-DD"
